$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The existing row 4 (MENINGITIS / 31663_7#1000 / 1234STDY1234 / SAMPLE / IB / Wellcome Sanger
# Institute) is being duplicated down into row 5 (without shifting the rows below, since rows
# 5 and 6 are currently blank/unused), and row 5 is then edited so that the sample is now keyed
# by sampleId instead of laneId: the laneId (B) is cleared, the sampleId (C) becomes a new id,
# and the publicName (D) becomes "NOLANEID".
$ws.Range("A4:F4").Copy($ws.Range("A5:F5"))

$ws.Cells.Item(5, 2).Value = ""
$ws.Cells.Item(5, 3).Value = "1234STDY1235"
$ws.Cells.Item(5, 4).Value = "NOLANEID"

# Row 4 itself is overwritten with padded/whitespace-dirty versions of the original values.
$ws.Cells.Item(4, 1).Value = "   MENINGITIS   "
$ws.Cells.Item(4, 2).Value = "   31663_7#1000   "
$ws.Cells.Item(4, 3).Value = "   1234STDY1234   "
$ws.Cells.Item(4, 4).Value = "   SAMPLE   "
$ws.Cells.Item(4, 5).Value = "   IB   "
$ws.Cells.Item(4, 6).Value = "       Wellcome Sanger Institute       "

$ws.Range("E24").Select()
